# Updates cached numeric values in the profit-tracking sheets to reflect
# refreshed market-board prices pulled by the scheduled runner.
# (No formulas are involved -- these are literal scraped values.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1667.5
$ws.Range("I15").Value = 1667.5
$ws.Range("K15").Value = 5002.5
$ws.Range("M15").Value = -4833.5
# Row 28
$ws.Range("H28").Value = 956.93335
$ws.Range("I28").Value = 911.8461
$ws.Range("J28").Value = 1250
$ws.Range("K28").Value = 911.8461
$ws.Range("L28").Value = 1250
$ws.Range("M28").Value = -426.8461
$ws.Range("N28").Value = -2220
# Row 107
$ws.Range("H107").Value = 586.6667
$ws.Range("I107").Value = 321.26666
$ws.Range("J107").Value = 1029
$ws.Range("K107").Value = 321.26666
$ws.Range("L107").Value = 1029
$ws.Range("M107").Value = 1598.73334
$ws.Range("N107").Value = -4869
# Row 132
$ws.Range("H132").Value = 1461.909
$ws.Range("I132").Value = 848.2059
$ws.Range("K132").Value = 2544.6177
$ws.Range("M132").Value = -14.61770000000024
# Row 137
$ws.Range("H137").Value = 872.60974
$ws.Range("I137").Value = 628.6818
$ws.Range("J137").Value = 1155.0526
$ws.Range("K137").Value = 1886.0454
$ws.Range("L137").Value = 3465.1578
$ws.Range("M137").Value = 663.9546
$ws.Range("N137").Value = -8565.157800000001
# Row 138
$ws.Range("H138").Value = 3097.0344
$ws.Range("I138").Value = 2831.125
$ws.Range("J138").Value = 4373.4
$ws.Range("K138").Value = 8493.375
$ws.Range("L138").Value = 13120.2
$ws.Range("M138").Value = -3353.375
$ws.Range("N138").Value = -23400.2

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1337.1818
$ws.Range("I45").Value = 1194.8948
$ws.Range("J45").Value = 1530.2858
$ws.Range("K45").Value = 1194.8948
$ws.Range("L45").Value = 1530.2858
$ws.Range("M45").Value = -817.8948
$ws.Range("N45").Value = -2284.2858
# Row 74
$ws.Range("H74").Value = 869.1923
$ws.Range("I74").Value = 439.58823
$ws.Range("J74").Value = 1680.6666
$ws.Range("K74").Value = 439.58823
$ws.Range("L74").Value = 1680.6666
$ws.Range("M74").Value = 434.41177
$ws.Range("N74").Value = -3428.6666
# Row 77
$ws.Range("H77").Value = 869.1923
$ws.Range("I77").Value = 439.58823
$ws.Range("J77").Value = 1680.6666
$ws.Range("K77").Value = 2197.94115
$ws.Range("L77").Value = 8403.333000000001
$ws.Range("M77").Value = 2170.05885
$ws.Range("N77").Value = -17139.333
# Row 122
$ws.Range("H122").Value = 1884.4857
$ws.Range("I122").Value = 1867.9259
$ws.Range("J122").Value = 1940.375
$ws.Range("K122").Value = 5603.7777
$ws.Range("L122").Value = 5821.125
$ws.Range("M122").Value = -3153.7777
$ws.Range("N122").Value = -10721.125

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 4317.6
$ws.Range("I99").Value = 5022
$ws.Range("K99").Value = 5022
$ws.Range("M99").Value = -3524

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 6976666.5
$ws.Range("J4").Value = 6976666.5
$ws.Range("L4").Value = 6976666.5
$ws.Range("N4").Value = -6976890.5
# Row 58
$ws.Range("H58").Value = 2668.6667
$ws.Range("I58").Value = 2456
$ws.Range("J58").Value = 2775
$ws.Range("K58").Value = 2456
$ws.Range("L58").Value = 2775
$ws.Range("M58").Value = -2253
$ws.Range("N58").Value = -3181
# Row 62
$ws.Range("H62").Value = 15922346
$ws.Range("I62").Value = 83583816
$ws.Range("J62").Value = 1999.1765
$ws.Range("K62").Value = 83583816
$ws.Range("L62").Value = 1999.1765
$ws.Range("M62").Value = -83583192
$ws.Range("N62").Value = -3247.1765
# Row 65
$ws.Range("H65").Value = 15922346
$ws.Range("I65").Value = 83583816
$ws.Range("J65").Value = 1999.1765
$ws.Range("K65").Value = 417919080
$ws.Range("L65").Value = 9995.8825
$ws.Range("M65").Value = -417915960
$ws.Range("N65").Value = -16235.8825
# Row 132
$ws.Range("H132").Value = 2413.75
$ws.Range("I132").Value = 2178.8215
$ws.Range("J132").Value = 3236
$ws.Range("K132").Value = 6536.4645
$ws.Range("L132").Value = 9708
$ws.Range("M132").Value = -4006.4645
$ws.Range("N132").Value = -14768
# Row 136
$ws.Range("H136").Value = 2668.6667
$ws.Range("I136").Value = 2456
$ws.Range("J136").Value = 2775
$ws.Range("K136").Value = 7368
$ws.Range("L136").Value = 8325
$ws.Range("M136").Value = -4818
$ws.Range("N136").Value = -13425

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1212.6522
$ws.Range("I5").Value = 576.9091
$ws.Range("J5").Value = 1795.4166
$ws.Range("K5").Value = 1730.7273
$ws.Range("L5").Value = 5386.2498
$ws.Range("M5").Value = -1618.7273
$ws.Range("N5").Value = -5610.2498
# Row 68
$ws.Range("H68").Value = 911.89
$ws.Range("I68").Value = 613.39685
$ws.Range("J68").Value = 1420.1351
$ws.Range("K68").Value = 1840.19055
$ws.Range("L68").Value = 4260.4053
$ws.Range("M68").Value = -1029.19055
$ws.Range("N68").Value = -5882.4053
# Row 71
$ws.Range("H71").Value = 911.89
$ws.Range("I71").Value = 613.39685
$ws.Range("J71").Value = 1420.1351
$ws.Range("K71").Value = 5520.57165
$ws.Range("L71").Value = 12781.2159
$ws.Range("M71").Value = -1464.57165
$ws.Range("N71").Value = -20893.2159
# Row 135
$ws.Range("H135").Value = 1212.6522
$ws.Range("I135").Value = 576.9091
$ws.Range("J135").Value = 1795.4166
$ws.Range("K135").Value = 5192.1819
$ws.Range("L135").Value = 16158.7494
$ws.Range("M135").Value = -2657.1819
$ws.Range("N135").Value = -21228.7494

$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 30500
$ws.Range("J6").Value = 30500
$ws.Range("L6").Value = 30500
$ws.Range("N6").Value = -30726
# Row 16
$ws.Range("H16").Value = 30500
$ws.Range("J16").Value = 30500
$ws.Range("L16").Value = 30500
$ws.Range("N16").Value = -31000
# Row 102
$ws.Range("H102").Value = 3325.4583
$ws.Range("I102").Value = 2665.1428
$ws.Range("J102").Value = 4249.9
$ws.Range("K102").Value = 2665.1428
$ws.Range("L102").Value = 4249.9
$ws.Range("M102").Value = -1043.1428
$ws.Range("N102").Value = -7493.9
# Row 122
$ws.Range("H122").Value = 4753.8335
$ws.Range("I122").Value = 4126.75
$ws.Range("J122").Value = 6008
$ws.Range("K122").Value = 12380.25
$ws.Range("L122").Value = 18024
$ws.Range("M122").Value = -9930.25
$ws.Range("N122").Value = -22924

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2546.484
$ws.Range("I40").Value = 2381.3333
$ws.Range("J40").Value = 2893.3
$ws.Range("K40").Value = 2893.3
$ws.Range("L40").Value = 2893.3
$ws.Range("M40").Value = -2245.3333
$ws.Range("N40").Value = -3165.3
